# Fruta / hortaliza, semanal
# Inserts a new weekly data row for "Zapallo" (Camote) above the existing
# row 51, shifting the remaining records (old rows 51-83) down by one
# (new rows 52-84), and populates the freshly inserted row 51 with the
# new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 51 (and everything below it) down by one row.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with this week's record.
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 44438
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 100112045
$ws.Cells.Item(51, 7).Value = "Zapallo"
$ws.Cells.Item(51, 8).Value = "Camote"
$ws.Cells.Item(51, 9).Value = "1a (guarda)"
$ws.Cells.Item(51, 10).Value = 300
$ws.Cells.Item(51, 11).Value = 600
$ws.Cells.Item(51, 12).Value = 650
$ws.Cells.Item(51, 13).Value = 625
$ws.Cells.Item(51, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(51, 15).Value = "Región del Maule"
$ws.Cells.Item(51, 16).Value = 625
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
